$d = $word.ActiveDocument

$replacements = @(
    @{old="83×91=7553"; new="40×11=440"},
    @{old="57×27=1539"; new="83×28=2324"},
    @{old="61×78=4758"; new="93×25=2325"},
    @{old="33×95=3135"; new="41×14=574"},
    @{old="50×21=1050"; new="66×55=3630"},
    @{old="87×14=1218"; new="83×23=1909"},
    @{old="90×52=4680"; new="21×22=462"},
    @{old="27×32=864";  new="35×21=735"},
    @{old="71×37=2627"; new="37×14=518"},
    @{old="29×69=2001"; new="88×79=6952"},
    @{old="17×62=1054"; new="63×33=2079"},
    @{old="63×94=5922"; new="97×98=9506"},
    @{old="31×86=2666"; new="21×19=399"},
    @{old="78×40=3120"; new="56×43=2408"},
    @{old="50×23=1150"; new="68×60=4080"},
    @{old="28×40=1120"; new="90×50=4500"},
    @{old="19×44=836";  new="33×87=2871"},
    @{old="99×19=1881"; new="22×77=1694"},
    @{old="12×80=960";  new="19×69=1311"},
    @{old="92×53=4876"; new="46×27=1242"},
    @{old="85×42=3570"; new="71×80=5680"},
    @{old="29×41=1189"; new="33×89=2937"},
    @{old="87×42=3654"; new="72×82=5904"},
    @{old="44×77=3388"; new="53×87=4611"},
    @{old="65×95=6175"; new="63×70=4410"}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2) | Out-Null
}
